$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-15 15:59:34"
$wsZh.Range("G2").Value = "2016-01-15 16:00:24"

# Sheet "de-de": update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-15 15:59:44"
$wsDe.Range("G2").Value = "2016-01-15 16:00:40"
